$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 418, shifting existing rows 418-447 down to 419-448
$ws.Rows.Item(418).Insert()

# Populate the new row 418 with the new data point
$ws.Cells.Item(418, 1).Value = 10
$ws.Cells.Item(418, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(418, 3).Value = "La Araucanía"
$ws.Cells.Item(418, 4).Value = 45021
$ws.Cells.Item(418, 5).Value = 9
$ws.Cells.Item(418, 6).Value = 100112017
$ws.Cells.Item(418, 7).Value = "Apio"
$ws.Cells.Item(418, 8).Value = "Americana (o)"
$ws.Cells.Item(418, 9).Value = "Primera"
$ws.Cells.Item(418, 10).Value = 95
$ws.Cells.Item(418, 11).Value = 11000
$ws.Cells.Item(418, 12).Value = 13000
$ws.Cells.Item(418, 13).Value = 12368
$ws.Cells.Item(418, 14).Value = "`$/docena de matas"
$ws.Cells.Item(418, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(418, 16).Value = 2061
$ws.Cells.Item(418, 17).Value = 6
$ws.Cells.Item(418, 18).Value = "Hortaliza"
